{"js": "// The document contains a \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\" (notification) bullet list. One of\n// the list items is a stand-alone paragraph whose entire text is \"\u0399.\u039a.\u03a5. \".\n// The edit removes that whole paragraph; the following list item\n// (\"${local_directorate}\") is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"\u0399.\u039a.\u03a5.\";\nconst target = paragraphs.items.find((p) => p.text && p.text.indexOf(marker) !== -1);\n\nif (!target) {\n  throw new Error('Could not find paragraph containing \"\u0399.\u039a.\u03a5.\"');\n}\n\ntarget.delete();\nawait context.sync();\n", "ps1": "# The document contains a \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\" (notification) bullet list. One of\n# the list items is a stand-alone paragraph whose entire text is \"\u0399.\u039a.\u03a5. \".\n# This edit removes that whole paragraph; the following list item\n# (\"${local_directorate}\") is left untouched.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\u0399.\u039a.\u03a5.\")\n\nif (-not $found) {\n    throw 'Could not find paragraph containing \"\u0399.\u039a.\u03a5.\"'\n}\n\n# Grow the hit to its enclosing paragraph (wdParagraph = 4) so the whole\n# paragraph - including its paragraph mark - is removed, exactly like the\n# diff deletes the entire <w:p>...</w:p> element.\n$rng.Expand(4) | Out-Null\n$rng.Delete()\n"}
